# Sprint6 Burndown workbook update
# Updates the underlying effort-tracking data for "Burndown Chart Sprint6".
# Dependent formulas (N, L columns, the Completed/Remaining/Ideal summary
# rows 24-26, and the burndown chart's cached series) recalculate
# automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart Sprint6")

$ws.Range("H6").Value = 0.25
$ws.Range("H7").Value = 0.25

$ws.Range("G10").Value = 0.25
$ws.Range("H10").Value = 0.25

$ws.Range("H11").Value = 0.25

$ws.Range("H12").Value = 0.25

$ws.Range("H13").Value = 0.25

$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0.25
$ws.Range("F14").Value = 0.75

$ws.Range("D15").Value = 1
$ws.Range("G15").Value = 0.5
$ws.Range("H15").Value = 0.5

# Update the view state to match where the sprint is now being reviewed.
$ws.Activate()
$ws.Range("L28:M28").Select()

$ws5 = $wb.Worksheets.Item("Burndown Chart Sprint5")
$ws5.Activate()
$ws5.Range("J18").Select()

$ws.Activate()
